# Slide 1, shape 186 ("TextBox 1") holds the "Created by ..." credit line:
#   "Created by Mariko Tagawa (marikotagawa@gmail.com), JICA volunteer"
# The edit:
#   1. Drops the Microsoft "keep hyperlink color" extension
#      (<a:extLst><ahyp:hlinkClr .../></a:extLst>) from the hyperlinked
#      "Mariko Tagawa" run, leaving a bare <a:hlinkClick r:id="rId3"/>.
#   2. Adds a trailing space to "Mariko Tagawa".
#   3. Removes the parenthesised email " (marikotagawa@gmail.com)",
#      leaving ", JICA volunteer" directly after the name.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(186)
$tr = $shp.TextFrame.TextRange

# "Mariko Tagawa" is characters 12-24 of the original run text.
$nameRun = $tr.Characters(12, 13)

# Re-assigning the hyperlink's Address (even to its own value) makes
# PowerPoint regenerate a plain <a:hlinkClick> without the hyperlink-color
# extension block.
$hyperlink = $nameRun.ActionSettings.Item(1).Hyperlink
$hyperlink.Address = $hyperlink.Address

# Give the name a trailing space: "Mariko Tagawa" -> "Mariko Tagawa ".
$nameRun.Text = "Mariko Tagawa "

# After the rename the text is:
#   "Created by Mariko Tagawa  (marikotagawa@gmail.com), JICA volunteer"
# Characters 26-50 are " (marikotagawa@gmail.com)" - delete them so the
# text reads "Created by Mariko Tagawa , JICA volunteer".
$toDelete = $tr.Characters(26, 25)
$toDelete.Text = ""
